$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$lo1 = $ws1.ListObjects.Item(1)
$lo2 = $ws2.ListObjects.Item(1)

# --- Sheet1 (Table1): add two new outbreak-path rows (Q28, Q29) ---
$newRow1 = $lo1.ListRows.Add()
$ws1.Range("A33").Value = 44390
$ws1.Range("A33").NumberFormat = $ws1.Range("A32").NumberFormat
$ws1.Range("B33").Value = "Q13 w"
$ws1.Range("C33").Value = "Q28"
$ws1.Range("D33").Value = "Queensland"
$ws1.Range("F33").Value = "Carindale Greek Community Centre"
$ws1.Range("G33").Value = "Alpha (B.1.1.7)"
$ws1.Range("H33").Value = "Isolated"

$newRow2 = $lo1.ListRows.Add()
$ws1.Range("A34").Value = 44390
$ws1.Range("A34").NumberFormat = $ws1.Range("A32").NumberFormat
$ws1.Range("B34").Value = "Q13 w"
$ws1.Range("C34").Value = "Q29"
$ws1.Range("D34").Value = "Queensland"
$ws1.Range("F34").Value = "Carindale Greek Community Centre"
$ws1.Range("G34").Value = "Alpha (B.1.1.7)"
$ws1.Range("H34").Value = "Isolated"

$ws1.Range("H34").Select()

# --- Sheet2 (Date_Colours): regenerate the colour-ramp column, drop the
#     "Diamond" shape (everything becomes "Circle") and append 3 new dates ---

# Add three new ListRows (dates 44397-44399) at the bottom of the table
$lo2.ListRows.Add() | Out-Null
$lo2.ListRows.Add() | Out-Null
$lo2.ListRows.Add() | Out-Null

# New, longer "Colour Code" gradient (column B), now covering 44367-44390
$colourCode = @(
    "#fbf8fb","#f7f1f7","#f3e9f4","#f0e2f0","#ecdbec","#e8d4e9","#e4cde5","#e0c6e1",
    "#dcbfdd","#d8b8da","#d4b1d6","#d0aad2","#cca3ce","#c89ccb","#c495c7","#c08ec3",
    "#bb87c0","#b780bc","#b37ab8","#af73b5","#ab6cb1","#a765ad","#a25eaa","#9e57a6"
)
for ($i = 0; $i -lt $colourCode.Length; $i++) {
    $ws2.Range("B" + (2 + $i)).Value = $colourCode[$i]
}
# Clear the old tail of column B that no longer has a colour-code value
$ws2.Range("B26:B31").ClearContents()

# "Colour Code dbRaevn" (column C): re-write explicitly row by row, now
# shifted down one row starting at row 27, with the new rows 32-34 filled in
$ws2.Range("C21:C24").ClearContents()
$ws2.Range("C25").Value = "#E8D1FF"
$ws2.Range("C26").ClearContents()
$ws2.Range("C27").Value = "#E8D1FF"
$ws2.Range("C28").Value = "#E2C5FF"
$ws2.Range("C29").Value = "#CC99FF"
$ws2.Range("C30").Value = "#9968FB"
$ws2.Range("C31").Value = "#FFF7DA"
$ws2.Range("C32").Value = "#FFECA0"
$ws2.Range("C33").Value = "#FADA67"
$ws2.Range("C34").Value = "#E4AC01"

# All shapes are now "Circle" (the "Diamond" rows at the bottom are gone)
$ws2.Range("D26:D31").Value = "Circle"

# Dates for the three newly-added rows
$ws2.Range("A32").Value = 44397
$ws2.Range("A32").NumberFormat = $ws2.Range("A31").NumberFormat
$ws2.Range("A33").Value = 44398
$ws2.Range("A33").NumberFormat = $ws2.Range("A31").NumberFormat
$ws2.Range("A34").Value = 44399
$ws2.Range("A34").NumberFormat = $ws2.Range("A31").NumberFormat
